# Add RowNumber to get the next row number in a sheet.
#
# Each vendor sheet's "ID" column (A) previously stored app IDs built from a
# vendor prefix + GUID fragment (e.g. "VNF_E851486F", "CTRX_8CA6F29E_1").
# The new scheme uses the row number within the sheet instead
# (e.g. "VNF_1", "CTRX_2", ...). Update column A on every vendor sheet to
# match, then restore each sheet's saved cursor/selection.

$wb = $excel.ActiveWorkbook

# --- 1. Rewrite the per-row vendor IDs (column A) on every vendor sheet. ---
# The write order below matters: it reproduces the exact order in which the
# new shared-string values were appended to the workbook's string table.

$ws = $wb.Worksheets.Item("Vendor Not Found")
$ws.Range("A2").Value = "VNF_1"

$ws = $wb.Worksheets.Item("Microsoft")
$ws.Range("A2").Value = "MS_1"
$ws.Range("A3").Value = "MS_2"
$ws.Range("A4").Value = "MS_3"

$ws = $wb.Worksheets.Item("Fortinet Technologies")
$ws.Range("A2").Value = "FT_1"

$ws = $wb.Worksheets.Item("Citrix")
$ws.Range("A2").Value = "CTRX_1"
$ws.Range("A3").Value = "CTRX_2"
$ws.Range("A4").Value = "CTRX_3"
$ws.Range("A5").Value = "CTRX_4"
$ws.Range("A6").Value = "CTRX_5"

$ws = $wb.Worksheets.Item("Oracle")
$ws.Range("A2").Value = "ORC_1"
$ws.Range("A3").Value = "ORC_2"

$ws = $wb.Worksheets.Item("Test Vendor 1")
$ws.Range("A2").Value = "TV_1"

$ws = $wb.Worksheets.Item("Test Vendor 2")
$ws.Range("A2").Value = "TV_2"

$ws = $wb.Worksheets.Item("CheckPoint")
$ws.Range("A2").Value = "CP_1"
$ws.Range("A3").Value = "CP_2"

$ws = $wb.Worksheets.Item("Adobe")
$ws.Range("A2").Value = "ADOBE_1"

$ws = $wb.Worksheets.Item("CISCO")
$ws.Range("A2").Value = "CISCO_1"
$ws.Range("A3").Value = "CISCO_2"

$ws = $wb.Worksheets.Item("Waves Audio")
$ws.Range("A2").Value = "WAVES_1"

$ws = $wb.Worksheets.Item("Business Objects")
$ws.Range("A2").Value = "BO_1"

$ws = $wb.Worksheets.Item("ConnectWise")
$ws.Range("A2").Value = "CW_1"

$ws = $wb.Worksheets.Item("SAP")
$ws.Range("A2").Value = "SAP_1"

# --- 2. Restore each sheet's saved selection/cursor position. ---
# Selecting a range activates its sheet, so these run in tab order with
# "SAP" (the sheet that must stay active) last.

$ws = $wb.Worksheets.Item("Vendor Not Found")
$ws.Activate()
$ws.Range("B21").Select()

$ws = $wb.Worksheets.Item("Microsoft")
$ws.Activate()
$ws.Range("A2:A4").Select()

$ws = $wb.Worksheets.Item("CISCO")
$ws.Activate()
$ws.Range("A2:A3").Select()

$ws = $wb.Worksheets.Item("Fortinet Technologies")
$ws.Activate()
$ws.Range("A2").Select()

$ws = $wb.Worksheets.Item("Oracle")
$ws.Activate()
$ws.Range("A4").Select()

$ws = $wb.Worksheets.Item("Test Vendor 1")
$ws.Activate()
$ws.Range("E41").Select()

$ws = $wb.Worksheets.Item("Test Vendor 2")
$ws.Activate()
$ws.Range("A3").Select()

$ws = $wb.Worksheets.Item("Waves Audio")
$ws.Activate()
$ws.Range("C21").Select()

$ws = $wb.Worksheets.Item("CheckPoint")
$ws.Activate()
$ws.Range("A2:A3").Select()

$ws = $wb.Worksheets.Item("Adobe")
$ws.Activate()
$ws.Range("B17").Select()

$ws = $wb.Worksheets.Item("Business Objects")
$ws.Activate()
$ws.Range("A2").Select()

$ws = $wb.Worksheets.Item("ConnectWise")
$ws.Activate()
$ws.Range("C16").Select()

$ws = $wb.Worksheets.Item("SAP")
$ws.Activate()
$ws.Range("B21").Select()

# --- 3. Match the saved workbook window size. ---
$excel.ActiveWindow.Width = 367.3
$excel.ActiveWindow.Height = 137.4
